$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- MSRP (column D) and model-year (column C) updates for existing rows ---
$ws.Range("D2").Value = 42220
$ws.Range("D3").Value = 46690
$ws.Range("D4").Value = 44910
$ws.Range("D5").Value = 48865
$ws.Range("D6").Value = 45150
$ws.Range("D7").Value = 49620
$ws.Range("D8").Value = 47315
$ws.Range("D9").Value = 51230
$ws.Range("D10").Value = 40000
$ws.Range("D11").Value = 45200
$ws.Range("D12").Value = 49000
$ws.Range("D13").Value = 41910
$ws.Range("D14").Value = 47110
$ws.Range("D15").Value = 50910
$ws.Range("D16").Value = 45800
$ws.Range("C29").Value = 2021
$ws.Range("D29").Value = 53250
$ws.Range("C30").Value = 2021
$ws.Range("D30").Value = 56340
$ws.Range("C31").Value = 2021
$ws.Range("D31").Value = 64515
$ws.Range("D32").Value = 86730
$ws.Range("D33").Value = 91730
$ws.Range("D35").Value = 37610
$ws.Range("D36").Value = 39010
$ws.Range("D37").Value = 39710
$ws.Range("D38").Value = 41110
$ws.Range("D39").Value = 44060
$ws.Range("D40").Value = 45460
$ws.Range("D41").Value = 40160
$ws.Range("D42").Value = 46610
$ws.Range("D43").Value = 46910
$ws.Range("C45").Value = 2021
$ws.Range("D45").Value = 76000
$ws.Range("C46").Value = 2021
$ws.Range("D46").Value = 79250
$ws.Range("C47").Value = 2021
$ws.Range("D47").Value = 79600
$ws.Range("C48").Value = 2021
$ws.Range("D48").Value = 82850
$ws.Range("D53").Value = 65975
$ws.Range("D54").Value = 97625
$ws.Range("D55").Value = 45170
$ws.Range("D56").Value = 46570
$ws.Range("D57").Value = 48000
$ws.Range("D58").Value = 49400
$ws.Range("D59").Value = 48650
$ws.Range("D60").Value = 50050
$ws.Range("D61").Value = 47820
$ws.Range("D62").Value = 51210
$ws.Range("D63").Value = 51300
$ws.Range("D70").Value = 93050
$ws.Range("D71").Value = 97610
$ws.Range("D72").Value = 101100
$ws.Range("D81").Value = 40000
$ws.Range("D82").Value = 45200
$ws.Range("D83").Value = 49000
$ws.Range("D84").Value = 45800
$ws.Range("D86").Value = 51050
$ws.Range("D87").Value = 52450
$ws.Range("D88").Value = 53620
$ws.Range("D89").Value = 54000
$ws.Range("D90").Value = 55400
$ws.Range("D91").Value = 57210
$ws.Range("D92").Value = 49335
$ws.Range("D93").Value = 50735
$ws.Range("D94").Value = 51985
$ws.Range("D95").Value = 48835
$ws.Range("D96").Value = 51010
$ws.Range("D97").Value = 51765
$ws.Range("D98").Value = 53230

# --- Row 43: drop the custom row height (45) that is no longer needed ---
$ws.Rows.Item(43).AutoFit()

# --- New row 99: UX 250h AWD Black Line Special Edition (MSRP coming soon) ---
$ws.Range("B99").Value = "UX 250h AWD BLACK LINE SPECIAL EDITION"
$ws.Range("C99").Value = 2021
$ws.Range("D73").Copy()
$ws.Range("D99").PasteSpecial(-4122)
$ws.Range("D99").Value = "COMING SOON"
$ws.Range("E2").Copy()
$ws.Range("E99").PasteSpecial(-4122)
$ws.Range("E99").Value = 1025

# --- Update the view: scroll position + active selection ---
$excel.ActiveWindow.ScrollRow = 73
$ws.Range("D97").Select()
